# Add a new "Si4362 Configuration" worksheet, positioned right before the
# "RP2040 Memory Map" sheet, containing the RX modem register-calculation
# table (commit: "Add more register reading for RX modem work").

$wb = $excel.ActiveWorkbook

$memMap = $wb.Worksheets.Item("RP2040 Memory Map")
$ws = $wb.Worksheets.Add($memMap)
$ws.Name = "Si4362 Configuration"

# --- Header row ------------------------------------------------------
$ws.Range("A1").Value = "Parameter"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Unit"
$ws.Range("D1").Value = "Note"
$ws.Range("D1").WrapText = $true

# --- Crystal oscillator + decimation chain ---------------------------
# (cell text is entered in this particular order so that new shared-string
# table entries line up with how the workbook was originally authored)
$ws.Range("A2").Value = "xosc_freq"
$ws.Range("B2").Value = 30000000
$ws.Range("C2").Value = "Hz"

$ws.Range("A3").Value = "ndec2"
$ws.Range("B3").Value = 4

$ws.Range("A4").Value = "ndec1"
$ws.Range("B4").Value = 5

$ws.Range("A5").Value = "ndec0"
$ws.Range("B5").Value = 1

$ws.Range("A11").Value = "rxosr"
$ws.Range("B11").Value = 12

$ws.Range("A6").Value = "dwn3byp"
$ws.Range("B6").Value = 0

$ws.Range("A7").Value = "dwn2byp"
$ws.Range("B7").Value = 0

# --- Sample rate -------------------------------------------------------
$ws.Range("A9").Value = "sample_rate"
$ws.Range("D9").Value = "Sample rate for receiver I/Q and Bit Clock Recovery (BCR)."
$ws.Range("D9").WrapText = $true
$ws.Range("C9").Value = "kbps"
$ws.Range("B9").Formula = "=B2/B3/B4/B5/IF(B6,1,3)/IF(B7,1,2)/1000"

# --- Oversampling / bit rate -------------------------------------------
$ws.Range("A12").Value = "rx_oversampling_rate"
$ws.Range("D11").Value = "12 bits, indicates 8x the desired oversampling rate."
$ws.Range("D11").WrapText = $true
$ws.Range("B12").Formula = "=B11/8"

$ws.Range("A14").Value = "rx_bit_rate"
$ws.Range("B14").Formula = "=B9/B12"

# --- Column widths -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 34.17
$ws.Columns.Item(4).ColumnWidth = 45.3

# --- View: freeze header row, zoom, active cell -------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 176
$ws.Range("D6").Select() | Out-Null
